$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price/volume text values are preserved as text (not
# auto-converted to numbers) by pre-formatting the Price/Volume columns as Text.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.898.87"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "1.879.19"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "0.7341"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").Value = "242.14"
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").Value = "0.3146"
$ws.Range("E8").Value = "  +0.33%  "
$ws.Range("D9").Value = "0.07134"
$ws.Range("E9").Value = "  +0.85%  "
$ws.Range("D10").Value = "24.43"
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("D11").Value = "0.08269"
$ws.Range("E11").Value = "  -1.64%  "
$ws.Range("D12").Value = "0.7513"
$ws.Range("E12").Value = "  +0.50%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.910.43"
$ws.Range("E13").Value = "  +2.38%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "5.322"
$ws.Range("E14").Value = "  -0.47%  "
$ws.Range("D15").Value = "92.47"
$ws.Range("E15").Value = "  +0.50%  "
$ws.Range("D16").Value = "29.936.97"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("D17").Value = "6.069"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").Value = "248.16"
$ws.Range("E18").Value = "  +3.12%  "
$ws.Range("D19").Value = "13.40"
$ws.Range("E19").Value = "  -0.94%  "
$ws.Range("D20").Value = "0.000007842"
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.175.23"
$ws.Range("E21").Value = "  +2.91%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("D23").Value = "0.9995"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "7.762"
$ws.Range("E24").Value = "  -2.09%  "
$ws.Range("D25").Value = "0.1545"
$ws.Range("E25").Value = "  -1.22%  "
$ws.Range("D26").Value = "9.192"
$ws.Range("E26").Value = "  -1.15%  "
$ws.Range("D27").Value = "163.30"
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").Value = "18.58"
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("D29").Value = "2.029"
$ws.Range("E29").Value = "  +0.46%  "
$ws.Range("D30").Value = "1.450"
$ws.Range("E30").Value = "  -1.78%  "
$ws.Range("D31").Value = "4.541"
$ws.Range("E31").Value = "  +0.37%  "
$ws.Range("D32").Value = "1.531"
$ws.Range("E32").Value = "  +0.27%  "
$ws.Range("D33").Value = "4.189"
$ws.Range("E33").Value = "  -2.36%  "
$ws.Range("D34").Value = "0.05291"
$ws.Range("E34").Value = "  -0.25%  "
$ws.Range("D35").Value = "1.239"
$ws.Range("E35").Value = "  +0.87%  "
$ws.Range("D36").Value = "0.7599"
$ws.Range("E36").Value = "  +1.55%  "
$ws.Range("D37").Value = "0.9999"
$ws.Range("D38").Value = "2.717"
$ws.Range("E38").Value = "  +0.71%  "
$ws.Range("D39").Value = "0.01938"
$ws.Range("E39").Value = "  -0.55%  "
$ws.Range("D40").Value = "2.756"
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("D41").Value = "0.4509"
$ws.Range("E41").Value = "  +1.16%  "
$ws.Range("D42").Value = "6.004"
$ws.Range("E42").Value = "  -1.19%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "0.8710"
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "71.64"
$ws.Range("E44").Value = "  -0.48%  "
$ws.Range("D45").Value = "1.064.54"
$ws.Range("E45").Value = "  -2.67%  "
$ws.Range("D46").Value = "104.82"
$ws.Range("E46").Value = "  +2.42%  "
$ws.Range("D47").Value = "1.002"
$ws.Range("E47").Value = "  +0.16%  "
$ws.Range("D48").Value = "1.838"
$ws.Range("E48").Value = "  +0.27%  "
$ws.Range("D49").Value = "7.500"
$ws.Range("E49").Value = "  -3.30%  "
$ws.Range("D50").Value = "2.061.02"
$ws.Range("E50").Value = "  +1.57%  "
$ws.Range("D51").Value = "2.870"
$ws.Range("E51").Value = "  -6.33%  "
